$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 282.2857
$ws.Range("I12").Value = 299.2
$ws.Range("J12").Value = 240
$ws.Range("K12").Value = 299.2
$ws.Range("L12").Value = 240
$ws.Range("M12").Value = -129.2
$ws.Range("N12").Value = -580

$ws.Range("H32").Value = 3514
$ws.Range("I32").Value = 2374
$ws.Range("J32").Value = 4414
$ws.Range("K32").Value = 2374
$ws.Range("L32").Value = 4414
$ws.Range("M32").Value = -2048
$ws.Range("N32").Value = -5066

$ws.Range("H33").Value = 2011.875
$ws.Range("I33").Value = 2139.8667
$ws.Range("J33").Value = 92
$ws.Range("K33").Value = 2139.8667
$ws.Range("L33").Value = 92
$ws.Range("M33").Value = -1910.8667
$ws.Range("N33").Value = -550

$ws.Range("H40").Value = 3283.1667
$ws.Range("I40").Value = 1674.75
$ws.Range("J40").Value = 6500
$ws.Range("K40").Value = 1674.75
$ws.Range("L40").Value = 6500
$ws.Range("M40").Value = -1499.75
$ws.Range("N40").Value = -6850

$ws.Range("H107").Value = 1787.6666
$ws.Range("I107").Value = 1775
$ws.Range("J107").Value = 1802.1428
$ws.Range("K107").Value = 1775
$ws.Range("L107").Value = 1802.1428
$ws.Range("M107").Value = 145
$ws.Range("N107").Value = -5642.1428

$ws.Range("H112").Value = 84012.91
$ws.Range("J112").Value = 102432.78
$ws.Range("L112").Value = 307298.34
$ws.Range("N112").Value = -309514.34

$ws.Range("H125").Value = 1940.3077
$ws.Range("I125").Value = 1177.75
$ws.Range("J125").Value = 2279.2222
$ws.Range("K125").Value = 10599.75
$ws.Range("L125").Value = 20512.9998
$ws.Range("M125").Value = -8139.75
$ws.Range("N125").Value = -25432.9998

$ws.Range("H129").Value = 28762.686
$ws.Range("I129").Value = 70582.38
$ws.Range("J129").Value = 4051.0454
$ws.Range("K129").Value = 211747.14
$ws.Range("L129").Value = 12153.1362
$ws.Range("M129").Value = -206747.14
$ws.Range("N129").Value = -22153.1362

$ws.Range("H131").Value = 1765.3334
$ws.Range("I131").Value = 1318.4
$ws.Range("K131").Value = 3955.2
$ws.Range("M131").Value = 1084.8

$ws.Range("H133").Value = 113900
$ws.Range("J133").Value = 113900
$ws.Range("L133").Value = 113900
$ws.Range("N133").Value = -124020

$ws.Range("H138").Value = 3290.4775
$ws.Range("I138").Value = 1208.0416
$ws.Range("K138").Value = 3624.1248
$ws.Range("M138").Value = 1515.8752

$ws.Range("H139").Value = 76259.60000000001
$ws.Range("J139").Value = 76259.60000000001
$ws.Range("L139").Value = 76259.60000000001
$ws.Range("N139").Value = -86539.60000000001

$ws.Range("H140").Value = 99282.5
$ws.Range("J140").Value = 99282.5
$ws.Range("L140").Value = 99282.5
$ws.Range("N140").Value = -109642.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2694.3914
$ws.Range("I4").Value = 2000
$ws.Range("K4").Value = 2000
$ws.Range("M4").Value = -1884

$ws.Range("H5").Value = 300.3125
$ws.Range("I5").Value = 307.55554
$ws.Range("K5").Value = 307.55554
$ws.Range("M5").Value = -195.55554

$ws.Range("H44").Value = 39034.668
$ws.Range("J44").Value = 39034.668
$ws.Range("L44").Value = 39034.668
$ws.Range("N44").Value = -40010.668

$ws.Range("H61").Value = 43602.082
$ws.Range("I61").Value = 5111.25
$ws.Range("K61").Value = 5111.25
$ws.Range("M61").Value = -4899.25

$ws.Range("H74").Value = 153262.58
$ws.Range("I74").Value = 102586.22
$ws.Range("K74").Value = 102586.22
$ws.Range("M74").Value = -101712.22

$ws.Range("H77").Value = 153262.58
$ws.Range("I77").Value = 102586.22
$ws.Range("K77").Value = 512931.1
$ws.Range("M77").Value = -508563.1

$ws.Range("H136").Value = 43602.082
$ws.Range("I136").Value = 5111.25
$ws.Range("K136").Value = 15333.75
$ws.Range("M136").Value = -12783.75

$ws.Range("H138").Value = 122428.25
$ws.Range("J138").Value = 122428.25
$ws.Range("L138").Value = 122428.25
$ws.Range("N138").Value = -132708.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300.3125
$ws.Range("I4").Value = 307.55554
$ws.Range("K4").Value = 307.55554
$ws.Range("M4").Value = -192.55554

$ws.Range("H20").Value = 2783.9302
$ws.Range("I20").Value = 2787.1667
$ws.Range("J20").Value = 2776.4614
$ws.Range("K20").Value = 2787.1667
$ws.Range("L20").Value = 2776.4614
$ws.Range("M20").Value = -2540.1667
$ws.Range("N20").Value = -3270.4614

$ws.Range("H22").Value = 236.125

$ws.Range("H86").Value = 130952.29
$ws.Range("I86").Value = 3228.3333
$ws.Range("J86").Value = 226745.25
$ws.Range("K86").Value = 3228.3333
$ws.Range("L86").Value = 226745.25
$ws.Range("M86").Value = -2105.3333
$ws.Range("N86").Value = -228991.25

$ws.Range("H89").Value = 130952.29
$ws.Range("I89").Value = 3228.3333
$ws.Range("J89").Value = 226745.25
$ws.Range("K89").Value = 16141.6665
$ws.Range("L89").Value = 1133726.25
$ws.Range("M89").Value = -10525.6665
$ws.Range("N89").Value = -1144958.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 68000
$ws.Range("J135").Value = 68000
$ws.Range("L135").Value = 68000
$ws.Range("N135").Value = -78140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1107.75
$ws.Range("J45").Value = 1343.6666
$ws.Range("L45").Value = 4030.9998
$ws.Range("N45").Value = -5094.9998

$ws.Range("H92").Value = 237.3077
$ws.Range("I92").Value = 261.33334
$ws.Range("J92").Value = 230.1
$ws.Range("K92").Value = 784.0000200000001
$ws.Range("L92").Value = 690.3
$ws.Range("M92").Value = 463.9999799999999
$ws.Range("N92").Value = -3186.3

$ws.Range("H94").Value = 2511.5
$ws.Range("I94").Value = 2511.5
$ws.Range("K94").Value = 7534.5
$ws.Range("M94").Value = -6858.5

$ws.Range("H122").Value = 76572.5
$ws.Range("J122").Value = 114202.625
$ws.Range("L122").Value = 1027823.625
$ws.Range("N122").Value = -1032723.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18341666
$ws.Range("J11").Value = 26683332
$ws.Range("L11").Value = 26683332
$ws.Range("N11").Value = -26683610

$ws.Range("H17").Value = 2498.3333
$ws.Range("I17").Value = 3100
$ws.Range("K17").Value = 3100
$ws.Range("M17").Value = -2932

$ws.Range("H70").Value = 6169
$ws.Range("I70").Value = 6254
$ws.Range("J70").Value = 5999
$ws.Range("K70").Value = 6254
$ws.Range("L70").Value = 5999
$ws.Range("M70").Value = -5984
$ws.Range("N70").Value = -6539

$ws.Range("H73").Value = 6169
$ws.Range("I73").Value = 6254
$ws.Range("J73").Value = 5999
$ws.Range("K73").Value = 6254
$ws.Range("L73").Value = 5999
$ws.Range("M73").Value = -5318
$ws.Range("N73").Value = -7871

$ws.Range("H113").Value = 2964.75
$ws.Range("I113").Value = 2974.1428
$ws.Range("K113").Value = 2974.1428
$ws.Range("M113").Value = -804.1428000000001

$ws.Range("H132").Value = 1216981.5
$ws.Range("I132").Value = 1486238.9
$ws.Range("K132").Value = 4458716.699999999
$ws.Range("M132").Value = -4456186.699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7368.4688
$ws.Range("I22").Value = 2579.2
$ws.Range("J22").Value = 8255.370000000001
$ws.Range("K22").Value = 2579.2
$ws.Range("L22").Value = 8255.370000000001
$ws.Range("M22").Value = -2284.2
$ws.Range("N22").Value = -8845.370000000001

$ws.Range("H27").Value = 7368.4688
$ws.Range("I27").Value = 2579.2
$ws.Range("J27").Value = 8255.370000000001
$ws.Range("K27").Value = 2579.2
$ws.Range("L27").Value = 8255.370000000001
$ws.Range("M27").Value = -2472.2
$ws.Range("N27").Value = -8469.370000000001

$ws.Range("H40").Value = 36549.734
$ws.Range("I40").Value = 41628.156
$ws.Range("K40").Value = 41628.156
$ws.Range("M40").Value = -41492.156

$ws.Range("H93").Value = 2813.7334
$ws.Range("I93").Value = 2671.077
$ws.Range("K93").Value = 2671.077
$ws.Range("M93").Value = -1423.077

$ws.Range("H132").Value = 2942.9795
$ws.Range("I132").Value = 2552.5476
$ws.Range("J132").Value = 5285.5713
$ws.Range("K132").Value = 7657.6428
$ws.Range("L132").Value = 15856.7139
$ws.Range("M132").Value = -5127.6428
$ws.Range("N132").Value = -20916.7139

$ws.Range("H136").Value = 5134.5713
$ws.Range("I136").Value = 3876.375
$ws.Range("K136").Value = 11629.125
$ws.Range("M136").Value = -9079.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4079.5454
$ws.Range("I14").Value = 3499
$ws.Range("J14").Value = 4208.5557
$ws.Range("K14").Value = 3499
$ws.Range("L14").Value = 4208.5557
$ws.Range("M14").Value = -3331
$ws.Range("N14").Value = -4544.5557

$ws.Range("H126").Value = 23789.666
$ws.Range("I126").Value = 25071.412
$ws.Range("K126").Value = 75214.236
$ws.Range("M126").Value = -72744.236

$ws.Range("H132").Value = 11208.76
$ws.Range("I132").Value = 8696.409
$ws.Range("K132").Value = 26089.227
$ws.Range("M132").Value = -23559.227

$ws.Range("H136").Value = 3470.7036
$ws.Range("I136").Value = 3335
$ws.Range("K136").Value = 10005
$ws.Range("M136").Value = -7455

$ws.Range("H141").Value = 113817.86
$ws.Range("J141").Value = 113817.86
$ws.Range("L141").Value = 113817.86
$ws.Range("N141").Value = -124177.86
